# Applies the "Add files via upload" / "Added to ppt" edits described in the
# commit diff to the active presentation.

$p = $ppt.ActivePresentation
$CR = [char]13

# ---------------------------------------------------------------------------
# 1) Footer date placeholder: 10/19/2020 -> 10/20/2020
#    (Slide Master + every Slide Layout that carries the cached date field.)
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tf = $shp.TextFrame
            if ($tf.HasText) {
                $txt = $tf.TextRange.Text
                if ($txt -eq "10/19/2020") {
                    $tf.TextRange.Text = "10/20/2020"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 3 - Team Members & Responsibilities: flesh out Hunter's role.
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3shape = $s3.Shapes.Item(8)
$s3tr = $s3shape.TextFrame.TextRange
$s3para = $s3tr.Paragraphs(1, 1)
$s3para.Text = "zzz-placeholder-zzz"
$s3para2 = $s3tr.Paragraphs(1, 1)
$s3para2.Text = "Hunter Hutchison " + [char]0x2013 + " Developer: Front-end and testing"

# ---------------------------------------------------------------------------
# 3) Slide 4 - Technology Stack & Toolchain: Pyaudio -> Sound Device.
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4shape = $s4.Shapes.Item(5)
$s4tr = $s4shape.TextFrame.TextRange
$s4para = $s4tr.Paragraphs(2, 1)
$s4para.Text = "zzz-placeholder-zzz"
$s4para2 = $s4tr.Paragraphs(2, 1)
$s4para2.Text = "The team will use Microsoft Visual Studio to code and implement the Sound Device library"

# ---------------------------------------------------------------------------
# 4) Slide 5 - Potential Future Milestones: mention Matlibplot.
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5shape = $s5.Shapes.Item(7)
$s5tr = $s5shape.TextFrame.TextRange
$s5para = $s5tr.Paragraphs(2, 1)
$s5para.Text = "zzz-placeholder-zzz"
$s5para2 = $s5tr.Paragraphs(2, 1)
$s5para2.Text = "Milestone 4 hopes to be able to visualize the WAV file using "
$s5para2.InsertAfter("Matlibplot")
$s5tr2 = $s5shape.TextFrame.TextRange
$s5para3 = $s5tr2.Paragraphs(2, 1)
$s5para3.InsertAfter(".")

# ---------------------------------------------------------------------------
# 5) Slide 6 - "Github Repository" -> "Current Progress" with new bullets.
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)

$s6title = $s6.Shapes.Item(3)
$s6titleTr = $s6title.TextFrame.TextRange
$s6titleTr.Text = "zzz-placeholder-zzz"
$s6titleTr2 = $s6title.TextFrame.TextRange
$s6titleTr2.Text = "Current Progress"

$s6body = $s6.Shapes.Item(4)
$s6bodyTr = $s6body.TextFrame.TextRange
$s6bodyTr.Text = (
    "Implementation of simple Python GUI allowing the addition of a record button and a play-back button" + $CR +
    "Currently we only have a pop-up application window that can spawn in other pop-ups to notify user of running functions" + $CR +
    "Plan to improve the GUI and be able to better adjust the record time"
)

# ---------------------------------------------------------------------------
# 6) Slide 7 - Communication & Workflow Plan: "teams chat" -> "team chat".
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$s7shape = $s7.Shapes.Item(5)
$s7tr = $s7shape.TextFrame.TextRange
$s7para = $s7tr.Paragraphs(1, 1)
$s7para.Text = "zzz-placeholder-zzz"
$s7para2 = $s7tr.Paragraphs(1, 1)
$s7para2.Text = "We are group number 3 and meet in the number 3 team chat "
